# Generate Report for Handback
# Adds two new handback entries (07f22b3b-... and ce44ac28-...) to the
# Overview / zh-cn / de-de sheets, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # Overview
$ws2 = $wb.Worksheets.Item(2)  # zh-cn
$ws3 = $wb.Worksheets.Item(3)  # de-de

$statusText = "Handed back: in sync with en-US"
$blue = 0xED9564   # OLE (BGR) form of RGB FF6495ED, matches the existing HyperLink font color
$dateFmt = "yyyy-mm-dd HH:mm:ss"

function Style-Hyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $blue
}

function Style-Date($range) {
    $range.NumberFormat = $dateFmt
}

# ---------------------------------------------------------------------
# Sheet "Overview": two new rows (6 and 7), columns A (hyperlinked file
# name), B and C (status text).
# ---------------------------------------------------------------------

$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/07f22b3b2bf34ca5b902c4c6baae3194/e2e/07f22b3b-2bf3-4ca5-b902-c4c6baae3194.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "07f22b3b-2bf3-4ca5-b902-c4c6baae3194.md") | Out-Null
Style-Hyperlink $ws1.Range("A6")
$ws1.Range("B6").Value = $statusText
$ws1.Range("C6").Value = $statusText

$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/ce44ac289d014602b293d67cbbaf5ed3/e2e/ce44ac28-9d01-4602-b293-d67cbbaf5ed3.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ce44ac28-9d01-4602-b293-d67cbbaf5ed3.md") | Out-Null
Style-Hyperlink $ws1.Range("A7")
$ws1.Range("B7").Value = $statusText
$ws1.Range("C7").Value = $statusText

# ---------------------------------------------------------------------
# Sheet "zh-cn": two new rows (6 and 7), columns A-H.
# ---------------------------------------------------------------------

# Row 6 - 07f22b3b-2bf3-4ca5-b902-c4c6baae3194
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/07f22b3b2bf34ca5b902c4c6baae3194/e2e/07f22b3b-2bf3-4ca5-b902-c4c6baae3194.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "07f22b3b-2bf3-4ca5-b902-c4c6baae3194.md") | Out-Null
Style-Hyperlink $ws2.Range("A6")
$ws2.Range("B6").Value = $statusText
$ws2.Hyperlinks.Add($ws2.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/07f22b3b2bf34ca5b902c4c6baae3194/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/07f22b3b-2bf3-4ca5-b902-c4c6baae3194.14254085f8f7cae4224c47ac3ba3ec67fd5a93b4.zh-cn.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "07f22b3b-2bf3-4ca5-b902-c4c6baae3194.14254085f8f7cae4224c47ac3ba3ec67fd5a93b4.zh-cn.xlf") | Out-Null
Style-Hyperlink $ws2.Range("C6")
$ws2.Range("D6").Value = "2016-02-17 04:54:03"
Style-Date $ws2.Range("D6")
$ws2.Hyperlinks.Add($ws2.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/07f22b3b2bf34ca5b902c4c6baae3194/e2e/07f22b3b-2bf3-4ca5-b902-c4c6baae3194.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "07f22b3b-2bf3-4ca5-b902-c4c6baae3194.md") | Out-Null
Style-Hyperlink $ws2.Range("E6")
$ws2.Hyperlinks.Add($ws2.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/07f22b3b2bf34ca5b902c4c6baae3194/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/07f22b3b-2bf3-4ca5-b902-c4c6baae3194.14254085f8f7cae4224c47ac3ba3ec67fd5a93b4.zh-cn.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "07f22b3b-2bf3-4ca5-b902-c4c6baae3194.14254085f8f7cae4224c47ac3ba3ec67fd5a93b4.zh-cn.xlf") | Out-Null
Style-Hyperlink $ws2.Range("F6")
$ws2.Range("G6").Value = "2016-02-17 04:54:49"
$ws2.Range("H6").Value = "Include"

# Row 7 - ce44ac28-9d01-4602-b293-d67cbbaf5ed3
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/ce44ac289d014602b293d67cbbaf5ed3/e2e/ce44ac28-9d01-4602-b293-d67cbbaf5ed3.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ce44ac28-9d01-4602-b293-d67cbbaf5ed3.md") | Out-Null
Style-Hyperlink $ws2.Range("A7")
$ws2.Range("B7").Value = $statusText
$ws2.Hyperlinks.Add($ws2.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ce44ac289d014602b293d67cbbaf5ed3/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ce44ac28-9d01-4602-b293-d67cbbaf5ed3.7214772345b417cebf6937dce591607708667748.zh-cn.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ce44ac28-9d01-4602-b293-d67cbbaf5ed3.7214772345b417cebf6937dce591607708667748.zh-cn.xlf") | Out-Null
Style-Hyperlink $ws2.Range("C7")
$ws2.Range("D7").Value = "2016-02-17 04:54:03"
Style-Date $ws2.Range("D7")
$ws2.Hyperlinks.Add($ws2.Range("E7"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ce44ac289d014602b293d67cbbaf5ed3/e2e/ce44ac28-9d01-4602-b293-d67cbbaf5ed3.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ce44ac28-9d01-4602-b293-d67cbbaf5ed3.md") | Out-Null
Style-Hyperlink $ws2.Range("E7")
$ws2.Hyperlinks.Add($ws2.Range("F7"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ce44ac289d014602b293d67cbbaf5ed3/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ce44ac28-9d01-4602-b293-d67cbbaf5ed3.7214772345b417cebf6937dce591607708667748.zh-cn.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ce44ac28-9d01-4602-b293-d67cbbaf5ed3.7214772345b417cebf6937dce591607708667748.zh-cn.xlf") | Out-Null
Style-Hyperlink $ws2.Range("F7")
$ws2.Range("G7").Value = "2016-02-17 04:54:49"
$ws2.Range("H7").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de": two new rows (6 and 7), columns A-H.
# ---------------------------------------------------------------------

# Row 6 - 07f22b3b-2bf3-4ca5-b902-c4c6baae3194
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/07f22b3b2bf34ca5b902c4c6baae3194/e2e/07f22b3b-2bf3-4ca5-b902-c4c6baae3194.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "07f22b3b-2bf3-4ca5-b902-c4c6baae3194.md") | Out-Null
Style-Hyperlink $ws3.Range("A6")
$ws3.Range("B6").Value = $statusText
$ws3.Hyperlinks.Add($ws3.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/07f22b3b2bf34ca5b902c4c6baae3194/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/07f22b3b-2bf3-4ca5-b902-c4c6baae3194.14254085f8f7cae4224c47ac3ba3ec67fd5a93b4.de-de.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "07f22b3b-2bf3-4ca5-b902-c4c6baae3194.14254085f8f7cae4224c47ac3ba3ec67fd5a93b4.de-de.xlf") | Out-Null
Style-Hyperlink $ws3.Range("C6")
$ws3.Range("D6").Value = "2016-02-17 04:54:13"
Style-Date $ws3.Range("D6")
$ws3.Hyperlinks.Add($ws3.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/07f22b3b2bf34ca5b902c4c6baae3194/e2e/07f22b3b-2bf3-4ca5-b902-c4c6baae3194.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "07f22b3b-2bf3-4ca5-b902-c4c6baae3194.md") | Out-Null
Style-Hyperlink $ws3.Range("E6")
$ws3.Hyperlinks.Add($ws3.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/07f22b3b2bf34ca5b902c4c6baae3194/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/07f22b3b-2bf3-4ca5-b902-c4c6baae3194.14254085f8f7cae4224c47ac3ba3ec67fd5a93b4.de-de.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "07f22b3b-2bf3-4ca5-b902-c4c6baae3194.14254085f8f7cae4224c47ac3ba3ec67fd5a93b4.de-de.xlf") | Out-Null
Style-Hyperlink $ws3.Range("F6")
$ws3.Range("G6").Value = "2016-02-17 04:55:07"
$ws3.Range("H6").Value = "Include"

# Row 7 - ce44ac28-9d01-4602-b293-d67cbbaf5ed3
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/ce44ac289d014602b293d67cbbaf5ed3/e2e/ce44ac28-9d01-4602-b293-d67cbbaf5ed3.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ce44ac28-9d01-4602-b293-d67cbbaf5ed3.md") | Out-Null
Style-Hyperlink $ws3.Range("A7")
$ws3.Range("B7").Value = $statusText
$ws3.Hyperlinks.Add($ws3.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ce44ac289d014602b293d67cbbaf5ed3/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ce44ac28-9d01-4602-b293-d67cbbaf5ed3.7214772345b417cebf6937dce591607708667748.de-de.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ce44ac28-9d01-4602-b293-d67cbbaf5ed3.7214772345b417cebf6937dce591607708667748.de-de.xlf") | Out-Null
Style-Hyperlink $ws3.Range("C7")
$ws3.Range("D7").Value = "2016-02-17 04:54:13"
Style-Date $ws3.Range("D7")
$ws3.Hyperlinks.Add($ws3.Range("E7"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ce44ac289d014602b293d67cbbaf5ed3/e2e/ce44ac28-9d01-4602-b293-d67cbbaf5ed3.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ce44ac28-9d01-4602-b293-d67cbbaf5ed3.md") | Out-Null
Style-Hyperlink $ws3.Range("E7")
$ws3.Hyperlinks.Add($ws3.Range("F7"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ce44ac289d014602b293d67cbbaf5ed3/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ce44ac28-9d01-4602-b293-d67cbbaf5ed3.7214772345b417cebf6937dce591607708667748.de-de.xlf", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ce44ac28-9d01-4602-b293-d67cbbaf5ed3.7214772345b417cebf6937dce591607708667748.de-de.xlf") | Out-Null
Style-Hyperlink $ws3.Range("F7")
$ws3.Range("G7").Value = "2016-02-17 04:55:07"
$ws3.Range("H7").Value = "Include"

Write-Host "Handback rows added."
